$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Mexican Liga MX - Tigres vs Pumas UNAM)
$ws.Range("H2").Value  = 6
$ws.Range("M2").Value  = 1.07
$ws.Range("N2").Value  = 3.55
$ws.Range("R2").Value  = 1.33
$ws.Range("S2").Value  = 3.45
$ws.Range("T2").Value  = 1.95
$ws.Range("U2").Value  = 1.86
$ws.Range("W2").Value  = 2.44
$ws.Range("X2").Value  = 15
$ws.Range("Y2").Value  = 21
$ws.Range("Z2").Value  = 55
$ws.Range("AA2").Value = 240
$ws.Range("AB2").Value = 8
$ws.Range("AC2").Value = 9.800000000000001
$ws.Range("AD2").Value = 26
$ws.Range("AE2").Value = 120
$ws.Range("AF2").Value = 9.800000000000001
$ws.Range("AG2").Value = 10.5
$ws.Range("AH2").Value = 24
$ws.Range("AI2").Value = 120
$ws.Range("AJ2").Value = 17
$ws.Range("AK2").Value = 19
$ws.Range("AL2").Value = 42
$ws.Range("AM2").Value = 170
$ws.Range("AN2").Value = 11
$ws.Range("AO2").Value = 170

# Row 3 (Mexican Liga MX - Toluca vs Santos Laguna)
$ws.Range("J3").Value  = 7.6
$ws.Range("N3").Value  = 8.199999999999999
$ws.Range("O3").Value  = 1.09
$ws.Range("R3").Value  = 2.16
$ws.Range("S3").Value  = 1.74
$ws.Range("T3").Value  = 1.77
$ws.Range("U3").Value  = 2.08
$ws.Range("W3").Value  = 4.8
$ws.Range("Y3").Value  = 75
$ws.Range("Z3").Value  = 190
$ws.Range("AA3").Value = 570
$ws.Range("AB3").Value = 17.5
$ws.Range("AC3").Value = 21
$ws.Range("AD3").Value = 55
$ws.Range("AE3").Value = 220
$ws.Range("AF3").Value = 12.5
$ws.Range("AG3").Value = 13.5
$ws.Range("AH3").Value = 30
$ws.Range("AI3").Value = 140
$ws.Range("AJ3").Value = 12
$ws.Range("AK3").Value = 14
$ws.Range("AL3").Value = 32
$ws.Range("AM3").Value = 120
$ws.Range("AN3").Value = 2.86
$ws.Range("AO3").Value = 170

# Row 5 (German Bundesliga - Augsburg vs Union Berlin)
$ws.Range("F5").Value  = 2.72
$ws.Range("H5").Value  = 2.92
$ws.Range("U5").Value  = 1.97
$ws.Range("AA5").Value = 1000
$ws.Range("AB5").Value = 9.6
$ws.Range("AE5").Value = 80
$ws.Range("AF5").Value = 17
$ws.Range("AJ5").Value = 44
$ws.Range("AN5").Value = 36
$ws.Range("AO5").Value = 42

# Row 6 (Italian Serie A - Como vs AC Milan)
$ws.Range("T6").Value  = 1.85
